# Build site at 2022-09-26 16:07:08 UTC
# Applies the LOM3255.xlsx content edit:
#  - Row "Objetivos:" description replaced with the professor's name
#  - A new "Docentes responsáveis:" row inserted
#  - "Programa resumido:" value changed to "Semestral"
#  - Long "Programa:" paragraph removed (cell now holds the activation date, per source)
#  - Long "Método:" paragraph removed (cell now holds the professor's name, per source)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Objetivos: replace the long description with the professor's name.
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

# 2. Insert a new row 12 for "Docentes responsáveis:" (pushes the old rows 12-20 to 13-21).
$ws.Rows.Item(12).Insert()
$ws.Range("A12").Value = "Docentes responsáveis:"

# 3. "Programa resumido:" (now row 13) value becomes "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# 4. "Programa:" (now row 15) loses its long paragraph; cell now shows the activation date.
$ws.Range("B15").Value = "15/07/2015"
$ws.Range("C15").Value = "15/07/2015"

# 5. "Método:" (now row 18) loses its long paragraph; cell now shows the professor's name again.
$ws.Range("B18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C18").Value = "519033 - Carlos Yujiro Shigue"

Write-Output "edit applied"
